# Generate Report for Handoff
#
# b.md has now been handed off for localization (zh-cn and de-de). Update
# the Overview sheet status for b.md, and update the per-locale detail
# sheets with the new handoff file name / handoff datetime, keeping the
# existing hyperlink (rId) but repointing its displayed text.

$wb = $excel.ActiveWorkbook

function Set-HyperlinkDisplay($ws, [string]$addr, [string]$text) {
    $links = @($ws.Hyperlinks)
    foreach ($hl in $links) {
        if ($hl.Range.Address() -eq $addr) {
            $hl.TextToDisplay = $text
        }
    }
}

# --- Overview sheet: b.md row (row 3) status -> "Ready for handoff" ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = "Ready for handoff"
$overview.Range("C3").Value = "Ready for handoff"

# --- zh-cn sheet: b.md row (row 3) ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("B3").Value = "Ready for handoff"
$zhcn.Range("C3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("D3").Value = "2016-02-06 03:54:30"
Set-HyperlinkDisplay $zhcn "`$C`$3" "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"

# --- de-de sheet: b.md row (row 3) ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("B3").Value = "Ready for handoff"
$dede.Range("C3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("D3").Value = "2016-02-06 03:54:41"
Set-HyperlinkDisplay $dede "`$C`$3" "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
